# live_trading_results.xlsx update
# Trade #105 (HighProbConvergence) closes via early_exit; a new MarketMaking
# trade #134 opens. Reflect this across the Summary, Strategy Status,
# All Trades, HighProbConvergence and MarketMaking sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.45               # Current Capital
$summary.Range("B4").Value = 0.5600000000000001    # Total P&L $
$summary.Range("B5").Value = 0.11                  # Total P&L %
$summary.Range("B6").Value = 105                   # Total Trades
$summary.Range("B7").Value = 51                    # Winning Trades
$summary.Range("B9").Value = 48.57                 # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status sheet - HighProbConvergence row (row 3)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C3").Value = 100.41
$status.Range("D3").Value = 12
$status.Range("E3").Value = 0.42
$status.Range("F3").Value = 0.41
$status.Range("G3").Value = 83.33

# ---------------------------------------------------------------------
# All Trades sheet - trade #105 (row 106) closes with early_exit
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Cells.Item(106, 7).Value = 0.93                    # Exit Price
$allTrades.Cells.Item(106, 8).Value = "CLOSED"                # Status
$allTrades.Cells.Item(106, 9).Value = 8.267300000000001       # P&L %
$allTrades.Cells.Item(106, 10).Value = 0.07000000000000001    # P&L $
$allTrades.Cells.Item(106, 11).Value = 100.41                 # Capital After
$allTrades.Cells.Item(106, 12).Value = "early_exit"            # Exit Reason
$allTrades.Cells.Item(106, 13).Value = 0.12                   # Duration (min)

# New trade #134 (MarketMaking) appended as row 135
$allTrades.Cells.Item(135, 1).Value = 134
$allTrades.Cells.Item(135, 2).Value = "'2026-02-18"
$allTrades.Cells.Item(135, 2).ClearFormats()
$allTrades.Cells.Item(135, 3).Value = "00:28:02"
$allTrades.Cells.Item(135, 4).Value = "MarketMaking"
$allTrades.Cells.Item(135, 5).Value = "DOWN"
$allTrades.Cells.Item(135, 6).Value = 0.858985
$allTrades.Cells.Item(135, 8).Value = "OPEN"
$allTrades.Cells.Item(135, 9).Value = 0
$allTrades.Cells.Item(135, 10).Value = 0
$allTrades.Cells.Item(135, 11).Value = 99.47967800952271
$allTrades.Cells.Item(135, 13).Value = 0
$allTrades.Cells.Item(135, 14).Value = 0
$allTrades.Cells.Item(135, 15).Value = 0
$allTrades.Cells.Item(135, 16).Value = 0.65
$allTrades.Cells.Item(135, 17).Value = "Wide spread capture: 392 bps vs avg 307 bps"

# ---------------------------------------------------------------------
# HighProbConvergence sheet - trade #105 (row 13) closes with early_exit
# ---------------------------------------------------------------------
$hpc = $wb.Worksheets.Item("HighProbConvergence")
$hpc.Cells.Item(13, 7).Value = 0.93                    # Exit Price
$hpc.Cells.Item(13, 8).Value = "CLOSED"                # Status
$hpc.Cells.Item(13, 9).Value = 8.267300000000001       # P&L %
$hpc.Cells.Item(13, 10).Value = 0.07000000000000001    # P&L $
$hpc.Cells.Item(13, 11).Value = 100.41                 # Capital After
$hpc.Cells.Item(13, 16).Value = "early_exit"            # Exit Reason
$hpc.Cells.Item(13, 17).Value = 0.12                   # Duration (min)

# ---------------------------------------------------------------------
# MarketMaking sheet - new trade #134 appended as row 55
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Cells.Item(55, 1).Value = 134
$mm.Cells.Item(55, 2).Value = "'2026-02-18"
$mm.Cells.Item(55, 2).ClearFormats()
$mm.Cells.Item(55, 3).Value = "00:28:02"
$mm.Cells.Item(55, 4).Value = "MarketMaking"
$mm.Cells.Item(55, 5).Value = "DOWN"
$mm.Cells.Item(55, 6).Value = 0.858985
$mm.Cells.Item(55, 8).Value = "OPEN"
$mm.Cells.Item(55, 9).Value = 0
$mm.Cells.Item(55, 10).Value = 0
$mm.Cells.Item(55, 11).Value = 99.47967800952271
$mm.Cells.Item(55, 12).Value = 0
$mm.Cells.Item(55, 13).Value = 0
$mm.Cells.Item(55, 14).Value = 0.65
$mm.Cells.Item(55, 15).Value = "Wide spread capture: 392 bps vs avg 307 bps"
$mm.Cells.Item(55, 17).Value = 0
